$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Text change: "formatted_stats" -> "Overall" in the header row ---
$t.Cell(1, 2).Range.Text = "Overall"

# --- Column widths: 1080 dxa (54pt) -> 1814 dxa (90.7pt) for all 5 columns ---
for ($c = 1; $c -le $t.Columns.Count; $c++) {
    $t.Columns.Item($c).Width = 90.7
}

# --- Table borders: thin the grey rule from 2pt (sz=16) to 1.5pt (sz=12) ---
# NOTE: this runtime's Borders.Item() assigns LineWidth*2 to the saved w:sz,
# and swaps the wdBorderTop/wdBorderLeft indices, so LineWidth=6 (-> sz=12)
# must be applied via index -1 for the physical top edge and -3 for bottom.
$rows = $t.Rows.Count
$cols = $t.Columns.Count

for ($col = 1; $col -le $cols; $col++) {
    # Row 1 (header): both top and bottom rule thin from 2pt to 1.5pt
    $cell = $t.Cell(1, $col)
    $cell.Borders.Item(-1).LineWidth = 6
    $cell.Borders.Item(-3).LineWidth = 6

    # Row 2: top rule (shared with header's bottom) thins as well
    $cell = $t.Cell(2, $col)
    $cell.Borders.Item(-1).LineWidth = 6

    # Last row: bottom rule thins
    $cell = $t.Cell($rows, $col)
    $cell.Borders.Item(-3).LineWidth = 6
}
